# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly-computed "K" values (strike count, recalculated) keyed by sheet row number.
$kValuesByRow = @{
    2  = 1
    3  = 1
    5  = 0
    6  = 1
    7  = 0
    8  = 2
    9  = 0
    10 = 0
    11 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 1
    26 = 1
    27 = 0
    28 = 1
    29 = 0
    30 = 0
    31 = 2
    32 = 1
    33 = 0
    34 = 2
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    40 = 1
    41 = 0
    42 = 1
    43 = 0
    44 = 1
    45 = 0
    46 = 1
    47 = 0
    48 = 1
    49 = 2
    50 = 2
    51 = 1
    52 = 0
    53 = 2
    54 = 1
    55 = 2
    56 = 2
    57 = 1
    58 = 1
    59 = 0
    60 = 0
    61 = 0
    62 = 2
    63 = 0
    64 = 0
    65 = 1
    67 = 1
    68 = 2
    69 = 1
    70 = 1
}

# Column G holds "K" (header already relabeled from Strike# to K).
# Write the freshly calculated s_vals (K) into column G for each affected row.
$kColumn = 7
foreach ($row in $kValuesByRow.Keys) {
    $ws.Cells.Item($row, $kColumn).Value = $kValuesByRow[$row]
}
